# registration and login done
#
# Applies:
#  - Students!A2:F2 row data is rearranged/updated (registration + login
#    fields consolidated; the old 6-column layout collapses to 5 used
#    columns and the Role column moves into D2, shifting what used to be
#    in D/E/F left by one; F2 is cleared).
#  - Students sheet becomes the active sheet/tab (workbook-level active
#    tab moves from Teachers to Students), with selection on Students
#    moving to E9.
#  - Teachers sheet no longer carries the tabSelected flag (consequence
#    of Students becoming the active/selected tab).

$wb = $excel.ActiveWorkbook

$students = $wb.Worksheets.Item("Students")

# --- Update Students row 2 data (registration + hashed login info) ---
$students.Range("A2").Value = "r230014@famt.ac.in"
$students.Range("B2").Value = "Samiya"
$students.Range("C2").Value = '$2b$12$J55QWV5Dai8.x34A7o9lwufpgAhn/4oWgfrFuPZtc2idJVDFJL25C'
$students.Range("D2").Value = "student"
$students.Range("E2").Value = "saamia.kb@gmail.com"
$students.Range("F2").Value = ""

# --- Make "Students" the active/selected sheet (was "Teachers") ---
$students.Activate()
$students.Range("E9").Select()
